# The deck's primary theme (ppt/theme/theme1.xml) is switched from the
# "Integral" color scheme to the stock "Office Theme" color scheme.
#
# PowerPoint's object model exposes the theme's 12 scheme colors as
# $Master.Theme.ThemeColorScheme.Item(1..12) in the fixed order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2,
#   5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
#   11 hlink, 12 folHlink
# Writing .RGB on each item rewrites the corresponding <a:srgbClr val="..."/>
# inside <a:clrScheme> of the theme part backing the slide master.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# dk1 = 000000
$colorScheme.Item(1).RGB = 0
# lt1 = FFFFFF
$colorScheme.Item(2).RGB = 16777215
# dk2 = 44546A
$colorScheme.Item(3).RGB = 6968388
# lt2 = E7E6E6
$colorScheme.Item(4).RGB = 15132391
# accent1 = 5B9BD5
$colorScheme.Item(5).RGB = 13998939
# accent2 = ED7D31
$colorScheme.Item(6).RGB = 3243501
# accent3 = A5A5A5
$colorScheme.Item(7).RGB = 10855845
# accent4 = FFC000
$colorScheme.Item(8).RGB = 49407
# accent5 = 4472C4
$colorScheme.Item(9).RGB = 12874308
# accent6 = 70AD47
$colorScheme.Item(10).RGB = 4697456
# hlink = 0563C1
$colorScheme.Item(11).RGB = 12673797
# folHlink = 954F72
$colorScheme.Item(12).RGB = 7491477
